$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Drilling Cost")

# ---------------------------------------------------------------------------
# Insert two new columns (E and F) before the existing "Arithmetic Return"
# columns, which slide from E/F/G to G/H/I.
# ---------------------------------------------------------------------------
$ws.Range("E1:F1").EntireColumn.Insert()

# ---------------------------------------------------------------------------
# New header row (row 3) text for the two inserted columns.
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = "Average Cost per year"
$ws.Range("F3").Value = "Arithmetic average return"

# ---------------------------------------------------------------------------
# Column E: average cost per year = SUM(B:D)/3, rows 4-51.
# (style created first so it lands on the same cellXfs slot the workbook
# ends up with)
# ---------------------------------------------------------------------------
$ws.Range("B5").Copy()
$ws.Range("E4:E51").PasteSpecial(-4122)
$ws.Range("E4:E51").Borders.LineStyle = -4142

for ($r = 4; $r -le 51; $r++) {
    $ws.Cells.Item($r, 5).Formula = "=SUM(B${r}:D${r})/3"
}

# Give F3 the same look as the other header cells (font/fill/border/wrap)
# plus a percentage number format.
$ws.Range("B3").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("F3").NumberFormat = "0.00%"
$ws.Range("F3").Value = "Arithmetic average return"

# ---------------------------------------------------------------------------
# Column F: arithmetic average return = year-over-year % change of column E.
# Row 4 has no prior year, so it stays blank (formatted only).
# ---------------------------------------------------------------------------
$ws.Range("B6").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("F4").Borders.LineStyle = -4142
$ws.Range("F4").NumberFormat = "0.00%"
$ws.Range("F4").Value = ""

$ws.Range("B5").Copy()
$ws.Range("F5:F51").PasteSpecial(-4122)
$ws.Range("F5:F51").Borders.LineStyle = -4142
$ws.Range("F5:F51").NumberFormat = "0.00%"

for ($r = 5; $r -le 51; $r++) {
    $prev = $r - 1
    $ws.Cells.Item($r, 6).Formula = "=(E${r}-E${prev})/E${prev}"
}

# ---------------------------------------------------------------------------
# Printable area tweak that came along with the edit.
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
